# Apply the "evolution diff date format YYYY, YYYYtQ and YYYY/MM" edits
# to the dataset sheet.
#
# Changes (row labels refer to the "name" column, col A):
#   - eco_energies (row 7):      start_date 2019 -> 2019/07 ; end_date 2024 -> 2024t2
#   - eau_potable (row 8):       nb_row 80000 -> 34444
#   - transport_pub (row 11):    nb_row 90000 -> 999990
#   - conso_energie (row 18):    start_date 2015 -> 2011
#   - revenus_menages (row 19):  end_date 2024 -> 2027
#
# Finally, leave the active selection on L18 (conso_energie / start_date),
# matching the last cell touched in the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# eco_energies: start_date / end_date now carry a finer-grained format
$ws.Range("L7").Value = "2019/07"
$ws.Range("M7").Value = "2024t2"

# eau_potable: corrected row count
$ws.Range("H8").Value = 34444

# transport_pub: corrected row count
$ws.Range("H11").Value = 999990

# revenus_menages: end_date updated
$ws.Range("M19").Value = "2027"

# conso_energie: start_date updated (edited last, so it is left selected)
$ws.Range("L18").Value = "2011"

# Match the workbook's last selection (L18) left by the edit
$ws.Range("L18").Select()
